$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.327.16'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.82%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '4.028.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.11%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.42%  '

# Row 6
$ws.Range('E6').Value = '  +5.21%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.703'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +14.43%  '

# Row 8
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.756'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.68%  '

# Row 10
$ws.Range('E10').Value = '  +0.70%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000327'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.13%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +13.19%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.81'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.57%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.673.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.75%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.020.89'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.85%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.59'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.14%  '

# Row 18
$ws.Range('E18').Value = '  -0.32%  '

# Row 19
$ws.Range('E19').Value = '  -0.79%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.111.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.81%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.57%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '99.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +12.63%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.93%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.73%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.32%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.00%  '

# Row 27
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.01%  '

# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +31.55%  '

# Row 29
$ws.Range('E29').Value = '  +2.27%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.53'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.11%  '

# Row 32
$ws.Range('E32').Value = '  +4.05%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '677.90'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.44%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.38%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '66.13'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.56%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.46%  '

# Row 37
$ws.Range('E37').Value = '  -2.28%  '

# Row 38
$ws.Range('E38').Value = '  +4.76%  '

# Row 39
$ws.Range('E39').Value = '  +11.50%  '

# Row 40
$ws.Range('E40').Value = '  -6.98%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.69%  '

# Row 42
$ws.Range('E42').Value = '  -0.12%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0491'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.41%  '

# Row 45
$ws.Range('E45').Value = '  +6.41%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.63'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.73%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.29%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.79%  '

# Row 49
$ws.Range('E49').Value = '  -5.01%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.57%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.76%  '
